$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (price) to stay text, matching the source data which
# uses localized/grouped number strings that Excel would otherwise
# auto-convert to numeric values.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.665.21"
$ws.Range("E2").Value = "  +8.51%  "

$ws.Range("D3").Value = "1.944.07"
$ws.Range("E3").Value = "  +7.04%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "341.67"
$ws.Range("E5").Value = "  +3.42%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "0.4779"
$ws.Range("E7").Value = "  +4.87%  "

$ws.Range("D8").Value = "0.4136"
$ws.Range("E8").Value = "  +8.74%  "

$ws.Range("D9").Value = "48.49"
$ws.Range("E9").Value = "  +5.61%  "

$ws.Range("D10").Value = "0.08248"
$ws.Range("E10").Value = "  +5.04%  "

$ws.Range("D11").Value = "1.041"
$ws.Range("E11").Value = "  +8.60%  "

$ws.Range("D12").Value = "22.68"
$ws.Range("E12").Value = "  +8.33%  "

$ws.Range("D13").Value = "1.945.26"
$ws.Range("E13").Value = "  +4.85%  "

$ws.Range("D14").Value = "6.185"
$ws.Range("E14").Value = "  +5.85%  "

$ws.Range("D15").Value = "7.416"
$ws.Range("E15").Value = "  +5.08%  "

$ws.Range("D16").Value = "92.42"
$ws.Range("E16").Value = "  +3.93%  "

$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("D18").Value = "0.00001065"
$ws.Range("E18").Value = "  +4.64%  "

$ws.Range("D19").Value = "0.06673"
$ws.Range("E19").Value = "  +1.40%  "

$ws.Range("D20").Value = "18.05"
$ws.Range("E20").Value = "  +5.61%  "

$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("D22").Value = "29.650.64"
$ws.Range("E22").Value = "  +8.43%  "

$ws.Range("E23").Value = "  +6.14%  "

$ws.Range("D24").Value = "11.25"
$ws.Range("E24").Value = "  +4.23%  "

$ws.Range("D25").Value = "2.284"
$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("D26").Value = "2.185.25"
$ws.Range("E26").Value = "  +5.48%  "

$ws.Range("D27").Value = "160.84"
$ws.Range("E27").Value = "  +3.36%  "

$ws.Range("E28").Value = "  +4.69%  "

$ws.Range("D29").Value = "2.193"
$ws.Range("E29").Value = "  +7.89%  "

$ws.Range("D30").Value = "5.629"
$ws.Range("E30").Value = "  +7.84%  "

$ws.Range("D31").Value = "122.29"
$ws.Range("E31").Value = "  +4.13%  "

$ws.Range("D32").Value = "1.026"
$ws.Range("E32").Value = "  +10.56%  "

$ws.Range("D33").Value = "0.09668"
$ws.Range("E33").Value = "  +4.07%  "

$ws.Range("D34").Value = "1.467"
$ws.Range("E34").Value = "  +12.46%  "

$ws.Range("D35").Value = "3.682"
$ws.Range("E35").Value = "  +3.19%  "

$ws.Range("D36").Value = "5.493"
$ws.Range("E36").Value = "  +5.49%  "

$ws.Range("D37").Value = "0.06288"
$ws.Range("E37").Value = "  +6.61%  "

$ws.Range("D38").Value = "0.02318"
$ws.Range("E38").Value = "  +6.62%  "

$ws.Range("D39").Value = "8.627"
$ws.Range("E39").Value = "  +7.31%  "

$ws.Range("D40").Value = "1.193"
$ws.Range("E40").Value = "  +4.88%  "

$ws.Range("D41").Value = "0.6110"
$ws.Range("E41").Value = "  +6.74%  "

$ws.Range("D42").Value = "10.71"
$ws.Range("E42").Value = "  +8.54%  "

$ws.Range("D43").Value = "0.1905"
$ws.Range("E43").Value = "  +5.13%  "

$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "1.272"
$ws.Range("E45").Value = "  -0.58%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5720"
$ws.Range("E46").Value = "  +6.43%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "12.54"
$ws.Range("E47").Value = "  +6.04%  "

$ws.Range("D48").Value = "2.342"
$ws.Range("E48").Value = "  +32.66%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.001"
$ws.Range("E49").Value = "  +7.34%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.07421"
$ws.Range("E50").Value = "  +13.06%  "

$ws.Range("D51").Value = "114.57"
$ws.Range("E51").Value = "  +4.77%  "
